$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1966.5
$ws.Range("I29").Value = 1966.5
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 5899.5
$ws.Range("L29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("N29").Value = -5618.5

$ws.Range("H38").Value = 3188.3157
$ws.Range("I38").Value = 121.38461
$ws.Range("J38").Value = 9833.333000000001
$ws.Range("K38").Value = 364.15383
$ws.Range("L38").Value = 29499.999
$ws.Range("M38").Value = 7.846170000000029
$ws.Range("N38").Value = -30243.999

$ws.Range("H40").Value = 15496.6
$ws.Range("I40").Value = 2612.25
$ws.Range("K40").Value = 2612.25
$ws.Range("M40").Value = -2437.25

$ws.Range("H43").Value = 12004.75
$ws.Range("I43").Value = 13898.6
$ws.Range("J43").Value = 8848.333000000001
$ws.Range("K43").Value = 13898.6
$ws.Range("L43").Value = 8848.333000000001
$ws.Range("M43").Value = -13829.6
$ws.Range("N43").Value = -8986.333000000001

$ws.Range("H113").Value = 4724.625
$ws.Range("J113").Value = 4912
$ws.Range("L113").Value = 4912
$ws.Range("N113").Value = -11420

$ws.Range("H116").Value = 4755.3335
$ws.Range("I116").Value = 4724.875
$ws.Range("K116").Value = 4724.875
$ws.Range("M116").Value = -1282.875

$ws.Range("H137").Value = 2518.8
$ws.Range("I137").Value = 1578.2
$ws.Range("K137").Value = 4734.6
$ws.Range("M137").Value = -2184.6

$ws.Range("H138").Value = 3088.7693
$ws.Range("I138").Value = 1891.0869
$ws.Range("J138").Value = 4038.6553
$ws.Range("K138").Value = 5673.2607
$ws.Range("L138").Value = 12115.9659
$ws.Range("M138").Value = -533.2606999999998
$ws.Range("N138").Value = -22395.9659

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 45000
$ws.Range("J44").Value = 45000
$ws.Range("L44").Value = 45000
$ws.Range("N44").Value = -45976

$ws.Range("H74").Value = 2099.697
$ws.Range("I74").Value = 1440.6072
$ws.Range("K74").Value = 1440.6072
$ws.Range("M74").Value = -566.6071999999999

$ws.Range("H77").Value = 2099.697
$ws.Range("I77").Value = 1440.6072
$ws.Range("K77").Value = 7203.036
$ws.Range("M77").Value = -2835.036

$ws.Range("H122").Value = 2337.3684
$ws.Range("I122").Value = 2077.5557
$ws.Range("K122").Value = 6232.6671
$ws.Range("M122").Value = -3782.6671

$ws.Range("H132").Value = 3639.5386
$ws.Range("I132").Value = 3668.3044
$ws.Range("K132").Value = 11004.9132
$ws.Range("M132").Value = -8474.913199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2194.6875
$ws.Range("I134").Value = 2207.3215
$ws.Range("J134").Value = 2106.25
$ws.Range("K134").Value = 6621.9645
$ws.Range("L134").Value = 6318.75
$ws.Range("M134").Value = -4086.9645
$ws.Range("N134").Value = -11388.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5055595
$ws.Range("I99").Value = 5560664.5
$ws.Range("K99").Value = 5560664.5
$ws.Range("M99").Value = -5559166.5

$ws.Range("H107").Value = 1158.3636
$ws.Range("J107").Value = 556
$ws.Range("L107").Value = 556
$ws.Range("N107").Value = -4396

$ws.Range("H126").Value = 5055595
$ws.Range("I126").Value = 5560664.5
$ws.Range("K126").Value = 16681993.5
$ws.Range("M126").Value = -16679523.5

$ws.Range("H134").Value = 4512.926
$ws.Range("I134").Value = 4561.4
$ws.Range("K134").Value = 13684.2
$ws.Range("M134").Value = -11149.2

$ws.Range("H141").Value = 79265.336
$ws.Range("J141").Value = 86859.60000000001
$ws.Range("L141").Value = 86859.60000000001
$ws.Range("N141").Value = -97219.60000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1636.4286
$ws.Range("I5").Value = 1954
$ws.Range("J5").Value = 842.5
$ws.Range("K5").Value = 5862
$ws.Range("L5").Value = 2527.5
$ws.Range("M5").Value = -5750
$ws.Range("N5").Value = -2751.5

$ws.Range("H38").Value = 142.83333
$ws.Range("I38").Value = 92.40000000000001
$ws.Range("J38").Value = 162.23077
$ws.Range("K38").Value = 277.2
$ws.Range("L38").Value = 486.69231
$ws.Range("M38").Value = 69.79999999999995
$ws.Range("N38").Value = -1180.69231

$ws.Range("H122").Value = 6533
$ws.Range("I122").Value = 766
$ws.Range("J122").Value = 12300
$ws.Range("K122").Value = 6894
$ws.Range("L122").Value = 110700
$ws.Range("M122").Value = -4444
$ws.Range("N122").Value = -115600

$ws.Range("H135").Value = 1636.4286
$ws.Range("I135").Value = 1954
$ws.Range("J135").Value = 842.5
$ws.Range("K135").Value = 17586
$ws.Range("L135").Value = 7582.5
$ws.Range("M135").Value = -15051
$ws.Range("N135").Value = -12652.5

$ws.Range("H140").Value = 1951.5
$ws.Range("I140").Value = 1411.9231
$ws.Range("K140").Value = 4235.7693
$ws.Range("M140").Value = 944.2307000000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 173.47058
$ws.Range("I2").Value = 213
$ws.Range("K2").Value = 213
$ws.Range("M2").Value = -100

$ws.Range("H43").Value = 2601.625
$ws.Range("I43").Value = 2601.625
$ws.Range("K43").Value = 2601.625
$ws.Range("M43").Value = -2450.625

$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1000
$ws.Range("K46").Value = 1000
$ws.Range("M46").Value = -844

$ws.Range("H57").Value = 13238.75
$ws.Range("I57").Value = 13238.75
$ws.Range("K57").Value = 13238.75
$ws.Range("M57").Value = -12418.75

$ws.Range("H80").Value = 14710.8
$ws.Range("I80").Value = 13736.4
$ws.Range("J80").Value = 15198
$ws.Range("K80").Value = 13736.4
$ws.Range("L80").Value = 15198
$ws.Range("M80").Value = -12738.4
$ws.Range("N80").Value = -17194

$ws.Range("H83").Value = 14710.8
$ws.Range("I83").Value = 13736.4
$ws.Range("J83").Value = 15198
$ws.Range("K83").Value = 68682
$ws.Range("L83").Value = 75990
$ws.Range("M83").Value = -63690
$ws.Range("N83").Value = -85974

$ws.Range("H105").Value = 52249.25
$ws.Range("J105").Value = 52249.25
$ws.Range("L105").Value = 52249.25
$ws.Range("N105").Value = -59237.25

$ws.Range("H107").Value = 400.6842
$ws.Range("I107").Value = 304.30768
$ws.Range("K107").Value = 304.30768
$ws.Range("M107").Value = 1615.69232

$ws.Range("H122").Value = 3176.0312
$ws.Range("I122").Value = 2250.8
$ws.Range("K122").Value = 6752.400000000001
$ws.Range("M122").Value = -4302.400000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 10947.462
$ws.Range("I40").Value = 7885.6875
$ws.Range("K40").Value = 7885.6875
$ws.Range("M40").Value = -7749.6875

$ws.Range("H46").Value = 2364.4546
$ws.Range("I46").Value = 2472.6667
$ws.Range("J46").Value = 2323.875
$ws.Range("K46").Value = 2472.6667
$ws.Range("L46").Value = 2323.875
$ws.Range("M46").Value = -2284.6667
$ws.Range("N46").Value = -2699.875

$ws.Range("H100").Value = 4527.1816
$ws.Range("I100").Value = 2759.8
$ws.Range("K100").Value = 2759.8
$ws.Range("M100").Value = -2218.8

$ws.Range("H132").Value = 4099.72
$ws.Range("I132").Value = 4207.579
$ws.Range("J132").Value = 3758.1667
$ws.Range("K132").Value = 12622.737
$ws.Range("L132").Value = 11274.5001
$ws.Range("M132").Value = -10092.737
$ws.Range("N132").Value = -16334.5001

$ws.Range("H136").Value = 1999.25
$ws.Range("I136").Value = 2070.5715
$ws.Range("J136").Value = 1500
$ws.Range("K136").Value = 6211.7145
$ws.Range("L136").Value = 4500
$ws.Range("M136").Value = -3661.7145
$ws.Range("N136").Value = -9600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15586

$ws.Range("H45").Value = 25525.182
$ws.Range("I45").Value = 58275
$ws.Range("J45").Value = 18247.445
$ws.Range("K45").Value = 58275
$ws.Range("L45").Value = 18247.445
$ws.Range("M45").Value = -57784
$ws.Range("N45").Value = -19229.445

$ws.Range("H105").Value = 46528.25
$ws.Range("J105").Value = 46528.25
$ws.Range("L105").Value = 46528.25
$ws.Range("N105").Value = -53516.25

$ws.Range("H127").Value = 74695
$ws.Range("I127").Value = 74695
$ws.Range("K127").Value = 74695
$ws.Range("M127").Value = -69735
